# Carjacking arrests by month YoY - add data for 2021-11-27 (through 11-19)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / title to reflect the new "through" date
$ws.Name = "Through 2021-11-19"

# Row 13 label
$ws.Range("A13").Value = "November (through 11-19)"

# Row 13 values
$ws.Range("C13").Value = 18
$ws.Range("D13").Value = 0.0526
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 39
$ws.Range("G13").Value = 0.093
$ws.Range("I13").Value = 76
$ws.Range("J13").Value = 0.0256
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 25
$ws.Range("P13").Value = 0.1667
$ws.Range("Q13").Value = 6
$ws.Range("R13").Value = 111
$ws.Range("S13").Value = 0.0513
$ws.Range("U13").Value = 130
$ws.Range("V13").Value = 0.0076

# Row 14 (Total) values
$ws.Range("C14").Value = 244
$ws.Range("D14").Value = 0.1191
$ws.Range("E14").Value = 56
$ws.Range("F14").Value = 473
$ws.Range("G14").Value = 0.1059
$ws.Range("I14").Value = 725
$ws.Range("J14").Value = 0.0799
$ws.Range("N14").Value = 53
$ws.Range("O14").Value = 459
$ws.Range("P14").Value = 0.1035
$ws.Range("Q14").Value = 60
$ws.Range("R14").Value = 1114
$ws.Range("S14").Value = 0.0511
$ws.Range("U14").Value = 1484
$ws.Range("V14").Value = 0.0566
